$wb = $excel.ActiveWorkbook

# --- Sheet "linear" ---
$wsLinear = $wb.Worksheets.Item("linear")
$wsLinear.Range("B2").Value = 0.004236178276831877
$wsLinear.Range("B3").Value = -0.01195584445922531
$wsLinear.Range("B4").Value = 1.35502927106283

# --- Sheet "non-linear" ---
$wsNonLinear = $wb.Worksheets.Item("non-linear")
$wsNonLinear.Range("B2").Value = 0.02995110748995026
$wsNonLinear.Range("B3").Value = 0.03429075493017843
$wsNonLinear.Range("B4").Value = 1.379513816544738
$wsNonLinear.Range("B5").Value = 0.003682384919700102
$wsNonLinear.Range("B6").Value = -0.02267663403262827
$wsNonLinear.Range("B7").Value = 1.333420091586157
